$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (left to right) ----
$ws.Range("A1").Value = "filtro"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "current_address"
$ws.Range("E1").Value = "permanent_address"

# ---- Data, written column by column (top to bottom) to match the
#      original shared-string insertion order ----
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

$ws.Range("B2").Value = "Alexis Castano "
$ws.Range("B3").Value = "Sara Montoya"
$ws.Range("B4").Value = "Manuela Restrepo"
$ws.Range("B5").Value = "Susana Quiroz"

$ws.Range("C2").Value = "AlexisCastano@gmail.com"
$ws.Range("C3").Value = "SaraMontoya@gmail.com"
$ws.Range("C4").Value = "ManuelaRestrepo@gmail.com"
$ws.Range("C5").Value = "SusanaQuiroz@gmail.com"

$ws.Range("D2").Value = "cll324 sur 10"
$ws.Range("D3").Value = "cll389 sur 11"
$ws.Range("D4").Value = "cll382 sur 12"
$ws.Range("D5").Value = "cll394 sur 13"

$ws.Range("E2").Value = "cll324 sur 10"
$ws.Range("E3").Value = "cll389 sur 12"
$ws.Range("E4").Value = "cll382 sur 13"
$ws.Range("E5").Value = "cll394 sur 14"

# ---- Hyperlinks on the email column ----
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:AlexisCastano@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:SaraMontoya@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:ManuelaRestrepo@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:SusanaQuiroz@gmail.com") | Out-Null

# ---- Column widths (best effort bestFit) ----
$ws.Columns.Item(2).ColumnWidth = 16.5
$ws.Columns.Item(3).ColumnWidth = 27.333333333333332
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 18.0

# ---- Selection ----
$ws.Range("E9").Select() | Out-Null
